$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 73, pushing existing rows 73-75 down to 75-77
$ws.Rows("73:74").Insert()

# New row 73 data
$ws.Cells.Item(73, 1).Value = 11
$ws.Cells.Item(73, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(73, 3).Value = "Bíobío"
$ws.Cells.Item(73, 4).Value = 44931
$ws.Cells.Item(73, 4).NumberFormat = $ws.Cells.Item(75, 4).NumberFormat
$ws.Cells.Item(73, 5).Value = 8
$ws.Cells.Item(73, 6).Value = "Fruta"
$ws.Cells.Item(73, 7).Value = 100103
$ws.Cells.Item(73, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(73, 9).Value = 100103002
$ws.Cells.Item(73, 10).Value = "Ciruela"
$ws.Cells.Item(73, 11).Value = "Black Amber"
$ws.Cells.Item(73, 12).Value = "Primera"
$ws.Cells.Item(73, 13).Value = 100
$ws.Cells.Item(73, 14).Value = 15000
$ws.Cells.Item(73, 15).Value = 16000
$ws.Cells.Item(73, 16).Value = 15500
$ws.Cells.Item(73, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(73, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(73, 19).Value = 861
$ws.Cells.Item(73, 20).Value = 18

# New row 74 data
$ws.Cells.Item(74, 1).Value = 11
$ws.Cells.Item(74, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(74, 3).Value = "Bíobío"
$ws.Cells.Item(74, 4).Value = 44931
$ws.Cells.Item(74, 4).NumberFormat = $ws.Cells.Item(75, 4).NumberFormat
$ws.Cells.Item(74, 5).Value = 8
$ws.Cells.Item(74, 6).Value = "Fruta"
$ws.Cells.Item(74, 7).Value = 100103
$ws.Cells.Item(74, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(74, 9).Value = 100103002
$ws.Cells.Item(74, 10).Value = "Ciruela"
$ws.Cells.Item(74, 11).Value = "Black Amber"
$ws.Cells.Item(74, 12).Value = "Segunda"
$ws.Cells.Item(74, 13).Value = 50
$ws.Cells.Item(74, 14).Value = 14000
$ws.Cells.Item(74, 15).Value = 14000
$ws.Cells.Item(74, 16).Value = 14000
$ws.Cells.Item(74, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(74, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(74, 19).Value = 778
$ws.Cells.Item(74, 20).Value = 18
